$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arrLeft = New-Object 'object[,]' 24,5
$arrLeft[0,0] = 1.02
$arrLeft[0,1] = 1.026685587568758
$arrLeft[0,2] = 1.034303725391316
$arrLeft[0,3] = 1.026883986014946
$arrLeft[0,4] = 1.041409459000215
$arrLeft[1,0] = 1.02
$arrLeft[1,1] = 1.02757554019607
$arrLeft[1,2] = 1.035119459966553
$arrLeft[1,3] = 1.027637310653201
$arrLeft[1,4] = 1.04238331349327
$arrLeft[2,0] = 1.02
$arrLeft[2,1] = 1.02815139956578
$arrLeft[2,2] = 1.035647660894482
$arrLeft[2,3] = 1.028125166039026
$arrLeft[2,4] = 1.043014163593116
$arrLeft[3,0] = 1.02
$arrLeft[3,1] = 1.028393489725466
$arrLeft[3,2] = 1.035869802697426
$arrLeft[3,3] = 1.02833035559902
$arrLeft[3,4] = 1.043279539351393
$arrLeft[4,0] = 1.02
$arrLeft[4,1] = 1.028434137626689
$arrLeft[4,2] = 1.035907106290292
$arrLeft[4,3] = 1.028364813363609
$arrLeft[4,4] = 1.043324106822342
$arrLeft[5,0] = 1.02
$arrLeft[5,1] = 1.02815463439231
$arrLeft[5,2] = 1.035650628827601
$arrLeft[5,3] = 1.028127907420043
$arrLeft[5,4] = 1.043017708905104
$arrLeft[6,0] = 1.02
$arrLeft[6,1] = 1.026986350375552
$arrLeft[6,2] = 1.03457932999597
$arrLeft[6,3] = 1.027138491025567
$arrLeft[6,4] = 1.041738431629005
$arrLeft[7,0] = 1.02
$arrLeft[7,1] = 1.024927742925553
$arrLeft[7,2] = 1.032694435016063
$arrLeft[7,3] = 1.025398175156111
$arrLeft[7,4] = 1.039489616055452
$arrLeft[8,0] = 1.02
$arrLeft[8,1] = 1.023555448611943
$arrLeft[8,2] = 1.031439856664543
$arrLeft[8,3] = 1.024240181747995
$arrLeft[8,4] = 1.037994145170531
$arrLeft[9,0] = 1.02
$arrLeft[9,1] = 1.022961270782418
$arrLeft[9,2] = 1.03089710736773
$arrLeft[9,3] = 1.023739302605656
$arrLeft[9,4] = 1.037347494333791
$arrLeft[10,0] = 1.02
$arrLeft[10,1] = 1.022740572994089
$arrLeft[10,2] = 1.030695581377734
$arrLeft[10,3] = 1.023553336363434
$arrLeft[10,4] = 1.03710743569788
$arrLeft[11,0] = 1.02
$arrLeft[11,1] = 1.022787913115475
$arrLeft[11,2] = 1.030738805957815
$arrLeft[11,3] = 1.023593222993117
$arrLeft[11,4] = 1.037158922901361
$arrLeft[12,0] = 1.02
$arrLeft[12,1] = 1.022943027699222
$arrLeft[12,2] = 1.030880447622873
$arrLeft[12,3] = 1.023723928884969
$arrLeft[12,4] = 1.037327648224125
$arrLeft[13,0] = 1.02
$arrLeft[13,1] = 1.023038599783256
$arrLeft[13,2] = 1.030967727743622
$arrLeft[13,3] = 1.023804472085134
$arrLeft[13,4] = 1.037431623563142
$arrLeft[14,0] = 1.02
$arrLeft[14,1] = 1.023594883062416
$arrLeft[14,2] = 1.031475887624044
$arrLeft[14,3] = 1.024273434944764
$arrLeft[14,4] = 1.038037080292524
$arrLeft[15,0] = 1.02
$arrLeft[15,1] = 1.023943834929244
$arrLeft[15,2] = 1.031794775471903
$arrLeft[15,3] = 1.024567748494729
$arrLeft[15,4] = 1.0384171086747
$arrLeft[16,0] = 1.02
$arrLeft[16,1] = 1.024147375920122
$arrLeft[16,2] = 1.031980824676309
$arrLeft[16,3] = 1.024739468577008
$arrLeft[16,4] = 1.038638859345402
$arrLeft[17,0] = 1.02
$arrLeft[17,1] = 1.024216778659105
$arrLeft[17,2] = 1.032044270629195
$arrLeft[17,3] = 1.02479802945641
$arrLeft[17,4] = 1.038714485231196
$arrLeft[18,0] = 1.02
$arrLeft[18,1] = 1.023906395350018
$arrLeft[18,2] = 1.031760556919018
$arrLeft[18,3] = 1.024536166050027
$arrLeft[18,4] = 1.038376326264114
$arrLeft[19,0] = 1.02
$arrLeft[19,1] = 1.022897350137306
$arrLeft[19,2] = 1.030838735593666
$arrLeft[19,3] = 1.023685436966177
$arrLeft[19,4] = 1.037277959049021
$arrLeft[20,0] = 1.02
$arrLeft[20,1] = 1.022262960433911
$arrLeft[20,2] = 1.030259586103733
$arrLeft[20,3] = 1.02315102830244
$arrLeft[20,4] = 1.036588161119585
$arrLeft[21,0] = 1.02
$arrLeft[21,1] = 1.022599258537026
$arrLeft[21,2] = 1.030566562385264
$arrLeft[21,3] = 1.023434282594542
$arrLeft[21,4] = 1.036953760795422
$arrLeft[22,0] = 1.02
$arrLeft[22,1] = 1.023923312674258
$arrLeft[22,2] = 1.031776018664073
$arrLeft[22,3] = 1.024550436636915
$arrLeft[22,4] = 1.038394753812698
$arrLeft[23,0] = 1.02
$arrLeft[23,1] = 1.025459928003733
$arrLeft[23,2] = 1.033181376250232
$arrLeft[23,3] = 1.025847704366928
$arrLeft[23,4] = 1.040070335922976
$ws.Range("B2:F25").Value = $arrLeft

$arrRight = New-Object 'object[,]' 24,6
$arrRight[0,0] = 1.03083600681697
$arrRight[0,1] = 1.031847650046887
$arrRight[0,2] = 1.037103585536791
$arrRight[0,3] = 1.029705330976416
$arrRight[0,4] = 1.044189055169015
$arrRight[0,5] = 1.014653901931322
$arrRight[1,0] = 1.030942097318963
$arrRight[1,1] = 1.032377774215887
$arrRight[1,2] = 1.037728462511716
$arrRight[1,3] = 1.030266410919716
$arrRight[1,4] = 1.044973102917291
$arrRight[1,5] = 1.01482900863025
$arrRight[2,0] = 1.031008520972645
$arrRight[2,1] = 1.032720178761042
$arrRight[2,2] = 1.038132516948554
$arrRight[2,3] = 1.030629209876251
$arrRight[2,4] = 1.045480518650982
$arrRight[2,5] = 1.01494209478125
$arrRight[3,0] = 1.031035912123558
$arrRight[3,1] = 1.032863975733304
$arrRight[3,2] = 1.038302312492397
$arrRight[3,3] = 1.03078166774318
$arrRight[3,4] = 1.045693855077817
$arrRight[3,5] = 1.014989583218723
$arrRight[4,0] = 1.031040479915175
$arrRight[4,1] = 1.032888111053327
$arrRight[4,2] = 1.038330817832093
$arrRight[4,3] = 1.030807262370808
$arrRight[4,4] = 1.045729676268566
$arrRight[4,5] = 1.014997553617577
$arrRight[5,0] = 1.031008889071436
$arrRight[5,1] = 1.03272210077371
$arrRight[5,2] = 1.038134786037399
$arrRight[5,3] = 1.030631247274096
$arrRight[5,4] = 1.045483369190398
$arrRight[5,5] = 1.014942729532773
$arrRight[6,0] = 1.030872320539151
$arrRight[6,1] = 1.03202693592135
$arrRight[6,2] = 1.037314823431779
$arrRight[6,3] = 1.029895003111619
$arrRight[6,4] = 1.04445400948468
$arrRight[6,5] = 1.014713125282653
$arrRight[7,0] = 1.030614676471285
$arrRight[7,1] = 1.030797259980668
$arrRight[7,2] = 1.035867833812215
$arrRight[7,3] = 1.028595729793558
$arrRight[7,4] = 1.042640852121151
$arrRight[7,5] = 1.014306869379214
$arrRight[8,0] = 1.030431542089097
$arrRight[8,1] = 1.029974379554545
$arrRight[8,2] = 1.034901827606655
$arrRight[8,3] = 1.027728329741638
$arrRight[8,4] = 1.041432633996729
$arrRight[8,5] = 1.014034937224588
$arrRight[9,0] = 1.030349556626752
$arrRight[9,1] = 1.029617342790616
$arrRight[9,2] = 1.034483231856987
$arrRight[9,3] = 1.027352461765031
$arrRight[9,4] = 1.040909609681034
$arrRight[9,5] = 1.013916932814513
$arrRight[10,0] = 1.03031870070685
$arrRight[10,1] = 1.029484615733733
$arrRight[10,2] = 1.03432770154521
$arrRight[10,3] = 1.02721280704055
$arrRight[10,4] = 1.040715357772547
$arrRight[10,5] = 1.013873062643962
$arrRight[11,0] = 1.030325337623804
$arrRight[11,1] = 1.029513090994578
$arrRight[11,2] = 1.034361065352377
$arrRight[11,3] = 1.027242765272907
$arrRight[11,4] = 1.040757024423001
$arrRight[11,5] = 1.0138824746615
$arrRight[12,0] = 1.030347014280407
$arrRight[12,1] = 1.029606373721614
$arrRight[12,2] = 1.034470376588506
$arrRight[12,3] = 1.027340918684009
$arrRight[12,4] = 1.040893552289177
$arrRight[12,5] = 1.013913307267122
$arrRight[13,0] = 1.030360316624679
$arrRight[13,1] = 1.029663834048692
$arrRight[13,2] = 1.034537720892448
$arrRight[13,3] = 1.027401388891421
$arrRight[13,4] = 1.040977674647156
$arrRight[13,5] = 1.013932299218647
$arrRight[14,0] = 1.030436926652101
$arrRight[14,1] = 1.029998059745204
$arrRight[14,2] = 1.034929602023807
$arrRight[14,3] = 1.027753269101772
$arrRight[14,4] = 1.041467348499739
$arrRight[14,5] = 1.014042763427799
$arrRight[15,0] = 1.030484263285086
$arrRight[15,1] = 1.030207517597223
$arrRight[15,2] = 1.035175336883937
$arrRight[15,3] = 1.027973920637654
$arrRight[15,4] = 1.04177454690783
$arrRight[15,5] = 1.014111986441427
$arrRight[16,0] = 1.030511614624001
$arrRight[16,1] = 1.030329620880289
$arrRight[16,2] = 1.035318640008742
$arrRight[16,3] = 1.028102596042925
$arrRight[16,4] = 1.041953744168827
$arrRight[16,5] = 1.014152338329103
$arrRight[17,0] = 1.030520896713405
$arrRight[17,1] = 1.03037124304856
$arrRight[17,2] = 1.035367497570434
$arrRight[17,3] = 1.028146466432675
$arrRight[17,4] = 1.042014848076879
$arrRight[17,5] = 1.014166093074577
$arrRight[18,0] = 1.0304792113288
$arrRight[18,1] = 1.03018505197816
$arrRight[18,2] = 1.035148974939323
$arrRight[18,3] = 1.027950249590392
$arrRight[18,4] = 1.041741586007283
$arrRight[18,5] = 1.014104562021838
$arrRight[19,0] = 1.030340642159824
$arrRight[19,1] = 1.029578907248052
$arrRight[19,2] = 1.034438188392742
$arrRight[19,3] = 1.027312016063039
$arrRight[19,4] = 1.040853347601218
$arrRight[19,5] = 1.013904228881351
$arrRight[20,0] = 1.030251187754386
$arrRight[20,1] = 1.029197177377558
$arrRight[20,2] = 1.03399102766657
$arrRight[20,3] = 1.026910498847446
$arrRight[20,4] = 1.040295008295313
$arrRight[20,5] = 1.013778051361274
$arrRight[21,0] = 1.030298829853238
$arrRight[21,1] = 1.029399598261875
$arrRight[21,2] = 1.034228100382081
$arrRight[21,3] = 1.027123372563399
$arrRight[21,4] = 1.04059098164787
$arrRight[21,5] = 1.013844961190471
$arrRight[22,0] = 1.030481494891702
$arrRight[22,1] = 1.030195203440578
$arrRight[22,2] = 1.035160886859721
$arrRight[22,3] = 1.027960945600968
$arrRight[22,4] = 1.041756479577448
$arrRight[22,5] = 1.014107916873898
$arrRight[23,0] = 1.030683291843074
$arrRight[23,1] = 1.031115711409474
$arrRight[23,2] = 1.036242157928973
$arrRight[23,3] = 1.028931842916661
$arrRight[23,4] = 1.043109504954755
$arrRight[23,5] = 1.014412090982241
$ws.Range("I2:N25").Value = $arrRight

Write-Output "Applied vm_pu updates to rows 2-25"
